$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric values in column B ---
$ws.Range("B3").Value = 3401
$ws.Range("B4").Value = 2164
$ws.Range("B5").Value = 64
$ws.Range("B6").Value = 196601
$ws.Range("B7").Value = 91
$ws.Range("B8").Value = 68637
$ws.Range("B9").Value = 31.72

# --- Update labels in column A ---
# Single-quoted strings avoid PowerShell trying to expand '$_' in "R$_..." text
$ws.Range("A9").Value = 'TM_EXAMES_LABORATORIAIS'
$ws.Range("A10").Value = 'R$_CAIXA_TOTAL'
$ws.Range("A11").Value = 'R$_EFETIVACAO_TOTAL'
$ws.Range("A12").Value = 'NOVOS_PACIENTES'

# --- Update text values in column B (rows 10-12) ---
# B10/B11 hold decimal-formatted monetary text and are naturally kept as text
# because they are not parseable as plain numbers.
$ws.Range("B10").Value = "122.648,70"
$ws.Range("B11").Value = "138.288,02"

# B12 holds "250", a value that looks like a plain integer, so a normal
# Value assignment would silently convert it to a number. Force it to stay
# text (as in the source file) by writing it through a Text-formatted cell
# and then pasting the original (default/general) formatting back over it.
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "250"
$ws.Range("A12").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
